# penambahan query di list hotel
# lihat yg kuning2 di list hotel
#
# Adds a "Detail summary order / TRX013" row to the "Transaksi" sheet
# (renaming the old "Detail Booking Order" row to "Detail summary order"
# and re-appending the original text under a new code), and inserts a
# "Saldo Booking Order Detail / BLNC004" row into the "saldo" sheet
# (renaming the existing row to "... - summary" and re-using the original
# text under the new code).

$wb = $excel.ActiveWorkbook

# --- Sheet "Transaksi" ---------------------------------------------------
$wsTransaksi = $wb.Worksheets.Item("Transaksi")

# Row 4 keeps its code (TRX011) but its label becomes "Detail summary order"
$wsTransaksi.Range("A4").Value = "Detail summary order"

# New row 6: the original "Detail Booking Order" label gets a new code (TRX013)
$wsTransaksi.Range("A6").Value = "Detail Booking Order"
$wsTransaksi.Range("B6").Value = "TRX013"

# --- Sheet "saldo" ---------------------------------------------------------
$wsSaldo = $wb.Worksheets.Item("saldo")

# Insert a new row at 4, shifting rows 4-5 down to 5-6
$wsSaldo.Rows.Item(4).Insert()

# Row 3 label becomes "... - summary" (same code, BLNC002)
$wsSaldo.Range("A3").Value = "Saldo Booking Order Detail - summary"

# New row 4: original "Saldo Booking Order Detail" label gets a new code (BLNC004)
$wsSaldo.Range("A4").Value = "Saldo Booking Order Detail"
$wsSaldo.Range("B4").Value = "BLNC004"

# --- Restore the on-screen selections seen in the edited workbook ---------
$wsTransaksi.Activate() | Out-Null
$wsTransaksi.Range("B7").Select() | Out-Null

$wsSaldo.Activate() | Out-Null
$wsSaldo.Range("C14").Select() | Out-Null
